# Edit matching the target diff:
#  1. Slide 6's table switches from the custom "Table_0" style to the
#     built-in table style {38BB4A7B-B99F-412D-87C8-596D2E8D5EE9}.
#  2. The presentation's (slide master's) theme colour palette changes
#     from the "Integral" scheme to the standard "Office" scheme.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{38BB4A7B-B99F-412D-87C8-596D2E8D5EE9}")
    }
}

# --- 2. Theme colours: Integral -> Office ---------------------------------
# RGB() values below use the standard COM colour encoding (0x00BBGGRR),
# i.e. R + G*256 + B*65536, matching PowerPoint's ColorFormat.RGB.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
